$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$startRow = 81
$endRow = 86

$rows = @(
    @("2026-02-01", "17:51:17", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:51:27", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:51:38", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:51:48", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:51:59", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:52:09", "17:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

# Column A holds date-like text ("2026-02-01"). Force Text format first so
# Excel doesn't auto-convert the literal into a date serial number, matching
# how the rest of the log stores these values as plain strings.
$colA = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))
$colA.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

# Restore the plain "Normal" style on the written range so the new rows keep
# the same (unstyled) appearance as every other row in the sheet.
$fullRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 6))
$fullRange.Style = "Normal"
